$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "67.231.13"
Set-TextValue "E2" "  -1.56%  "

# Row 3
Set-TextValue "D3" "3.752.35"
Set-TextValue "E3" "  -2.75%  "

# Row 4
Set-TextValue "E4" "  +0.04%  "

# Row 5
Set-TextValue "D5" "593.69"
Set-TextValue "E5" "  -1.30%  "

# Row 6
Set-TextValue "D6" "167.89"
Set-TextValue "E6" "  -0.05%  "

# Row 7
Set-TextValue "D7" "3.751.82"
Set-TextValue "E7" "  -2.80%  "

# Row 8
Set-TextValue "E8" "  -0.07%  "

# Row 9
Set-TextValue "D9" "0.526"
Set-TextValue "E9" "  -0.69%  "

# Row 10
Set-TextValue "D10" "0.162"
Set-TextValue "E10" "  -1.01%  "

# Row 11
Set-TextValue "D11" "6.43"
Set-TextValue "E11" "  -0.19%  "

# Row 12
Set-TextValue "D12" "0.453"
Set-TextValue "E12" "  -1.17%  "

# Row 13
Set-TextValue "D13" "0.0000266"
Set-TextValue "E13" "  +1.26%  "

# Row 14
Set-TextValue "D14" "36.29"
Set-TextValue "E14" "  -2.11%  "

# Row 15
Set-TextValue "D15" "4.383.53"
Set-TextValue "E15" "  -2.88%  "

# Row 16
Set-TextValue "D16" "3.752.00"
Set-TextValue "E16" "  -3.00%  "

# Row 17
Set-TextValue "D17" "18.80"
Set-TextValue "E17" "  +3.39%  "

# Row 18
Set-TextValue "D18" "67.208.41"
Set-TextValue "E18" "  -1.83%  "

# Row 19
Set-TextValue "D19" "7.18"
Set-TextValue "E19" "  -2.69%  "

# Row 20
Set-TextValue "E20" "  +0.53%  "

# Row 21
Set-TextValue "D21" "10.51"
Set-TextValue "E21" "  -4.22%  "

# Row 22
Set-TextValue "D22" "465.27"
Set-TextValue "E22" "  -0.20%  "

# Row 23
Set-TextValue "D23" "0.717"
Set-TextValue "E23" "  -2.53%  "

# Row 24 (was PEPE, now Litecoin)
Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "83.55"
Set-TextValue "E24" "  +0.63%  "

# Row 25 (was Litecoin, now PEPE)
Set-TextValue "B25" "PEPE"
Set-TextValue "C25" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D25" "0.0000145"
Set-TextValue "E25" "  -9.44%  "

# Row 26
Set-TextValue "D26" "2.20"
Set-TextValue "E26" "  -1.67%  "

# Row 27
Set-TextValue "D27" "12.04"
Set-TextValue "E27" "  -0.60%  "

# Row 28
Set-TextValue "D28" "10.15"
Set-TextValue "E28" "  +0.90%  "

# Row 29
Set-TextValue "E29" "  -0.14%  "

# Row 30
Set-TextValue "D30" "2.89"
Set-TextValue "E30" "  -2.65%  "

# Row 31
Set-TextValue "D31" "3.906.03"
Set-TextValue "E31" "  -2.70%  "

# Row 32
Set-TextValue "D32" "7.55"
Set-TextValue "E32" "  -0.78%  "

# Row 33
Set-TextValue "D33" "2.23"
Set-TextValue "E33" "  -3.61%  "

# Row 34
Set-TextValue "D34" "30.19"
Set-TextValue "E34" "  -3.31%  "

# Row 35
Set-TextValue "D35" "9.06"
Set-TextValue "E35" "  -5.37%  "

# Row 36
Set-TextValue "D36" "3.718.38"
Set-TextValue "E36" "  -2.80%  "

# Row 37
Set-TextValue "D37" "3.78"
Set-TextValue "E37" "  +3.82%  "

# Row 38
Set-TextValue "D38" "0.104"
Set-TextValue "E38" "  -0.63%  "

# Row 39 (was Filecoin, now Kaspa)
Set-TextValue "B39" "Kaspa"
Set-TextValue "C39" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.137"
Set-TextValue "E39" "  -2.77%  "

# Row 40 (was Kaspa, now Filecoin)
Set-TextValue "B40" "Filecoin"
Set-TextValue "C40" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D40" "5.84"
Set-TextValue "E40" "  -1.33%  "

# Row 41
Set-TextValue "D41" "0.994"
Set-TextValue "E41" "  -2.64%  "

# Row 42
Set-TextValue "E42" "  -0.18%  "

# Row 43
Set-TextValue "D43" "0.311"
Set-TextValue "E43" "  -0.74%  "

# Row 44
Set-TextValue "E44" "  -0.02%  "

# Row 45
Set-TextValue "D45" "8.63"
Set-TextValue "E45" "  -0.09%  "

# Row 46
Set-TextValue "D46" "1.93"
Set-TextValue "E46" "  -2.43%  "

# Row 47
Set-TextValue "D47" "45.36"
Set-TextValue "E47" "  -3.64%  "

# Row 48
Set-TextValue "D48" "396.75"
Set-TextValue "E48" "  -5.84%  "

# Row 49 (was FLOKI, now Monero)
Set-TextValue "B49" "Monero"
Set-TextValue "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D49" "141.47"
Set-TextValue "E49" "  -0.22%  "

# Row 50 (was Monero, now FLOKI)
Set-TextValue "B50" "FLOKI"
Set-TextValue "C50" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D50" "0.000265"
Set-TextValue "E50" "  -11.30%  "

# Row 51
Set-TextValue "D51" "0.0351"
Set-TextValue "E51" "  -1.96%  "
